# Auto-generated edit script applying numeric value changes described by the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 77
$ws.Range("I2").Value = 77
$ws.Range("M2").Value = 36
$ws.Range("K2").Value = 77
$ws.Range("K4").Value = 247
$ws.Range("I4").Value = 247
$ws.Range("M4").Value = -133
$ws.Range("H4").Value = 247
$ws.Range("H9").Value = 251.875
$ws.Range("I9").Value = 192.33333
$ws.Range("K9").Value = 192.33333
$ws.Range("M9").Value = -23.33332999999999
$ws.Range("M12").Value = 0.8000000000000114
$ws.Range("H12").Value = 19981.2
$ws.Range("K12").Value = 169.2
$ws.Range("I12").Value = 169.2
$ws.Range("L39").Value = 627
$ws.Range("N39").Value = -1219
$ws.Range("M39").Value = -42857110
$ws.Range("J39").Value = 209
$ws.Range("H39").Value = 12500103
$ws.Range("K39").Value = 42857406
$ws.Range("I39").Value = 14285802
$ws.Range("I132").Value = 4776.355
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 14329.065
$ws.Range("H132").Value = 4752.0938
$ws.Range("N132").Value = -17060
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -11799.065

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J8").Value = 100000
$ws.Range("L8").Value = 100000
$ws.Range("N8").Value = -100288
$ws.Range("H8").Value = 2550000
$ws.Range("J101").Value = 98656.336
$ws.Range("H101").Value = 98656.336
$ws.Range("L101").Value = 98656.336
$ws.Range("N101").Value = -105146.336
$ws.Range("H104").Value = 41701.332
$ws.Range("J104").Value = 41701.332
$ws.Range("L104").Value = 41701.332
$ws.Range("N104").Value = -48689.332
$ws.Range("I132").Value = 2217.6667
$ws.Range("J132").Value = 4063.5
$ws.Range("K132").Value = 6653.000100000001
$ws.Range("H132").Value = 2646.9302
$ws.Range("N132").Value = -17250.5
$ws.Range("L132").Value = 12190.5
$ws.Range("M132").Value = -4123.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J44").Value = 20050
$ws.Range("L44").Value = 20050
$ws.Range("H44").Value = 20050
$ws.Range("N44").Value = -21044
$ws.Range("J80").Value = 735.5333000000001
$ws.Range("H80").Value = 746.1579
$ws.Range("K80").Value = 786
$ws.Range("M80").Value = 212
$ws.Range("N80").Value = -2731.5333
$ws.Range("L80").Value = 735.5333000000001
$ws.Range("I80").Value = 786
$ws.Range("M83").Value = 1062
$ws.Range("I83").Value = 786
$ws.Range("N83").Value = -13661.6665
$ws.Range("L83").Value = 3677.6665
$ws.Range("J83").Value = 735.5333000000001
$ws.Range("K83").Value = 3930
$ws.Range("H83").Value = 746.1579
$ws.Range("J94").Value = 3500
$ws.Range("N94").Value = -4402
$ws.Range("K94").Value = 2623.375
$ws.Range("L94").Value = 3500
$ws.Range("M94").Value = -2172.375
$ws.Range("I94").Value = 2623.375
$ws.Range("H94").Value = 2720.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 71.09090999999999
$ws.Range("I7").Value = 73.25
$ws.Range("K7").Value = 73.25
$ws.Range("M7").Value = 39.75
$ws.Range("L22").Value = 190
$ws.Range("H22").Value = 159.75
$ws.Range("K22").Value = 149.66667
$ws.Range("N22").Value = -890
$ws.Range("J22").Value = 190
$ws.Range("M22").Value = 200.33333
$ws.Range("I22").Value = 149.66667
$ws.Range("L43").Value = 17500
$ws.Range("J43").Value = 17500
$ws.Range("H43").Value = 17500
$ws.Range("N43").Value = -17868
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("I51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("H54").Value = 28250
$ws.Range("L54").Value = 26500
$ws.Range("K54").Value = 30000
$ws.Range("N54").Value = -27816
$ws.Range("I54").Value = 30000
$ws.Range("J54").Value = 26500
$ws.Range("M54").Value = -29342
$ws.Range("M61").ClearContents()
$ws.Range("I61").Value = 0
$ws.Range("H61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("N64").Value = -71295.60000000001
$ws.Range("L64").Value = 70799.60000000001
$ws.Range("J64").Value = 70799.60000000001
$ws.Range("H64").Value = 56142.715
$ws.Range("H67").Value = 56142.715
$ws.Range("N67").Value = -72515.60000000001
$ws.Range("L67").Value = 70799.60000000001
$ws.Range("J67").Value = 70799.60000000001
$ws.Range("N68").Value = -56493
$ws.Range("L68").Value = 54995
$ws.Range("H68").Value = 54995
$ws.Range("J68").Value = 54995
$ws.Range("N71").Value = -172473
$ws.Range("J71").Value = 54995
$ws.Range("L71").Value = 164985
$ws.Range("H71").Value = 54995
$ws.Range("J101").Value = 17500
$ws.Range("H101").Value = 17500
$ws.Range("L101").Value = 17500
$ws.Range("N101").Value = -23990
$ws.Range("J110").Value = 69989.5
$ws.Range("L110").Value = 69989.5
$ws.Range("N110").Value = -78169.5
$ws.Range("H110").Value = 69989.5
$ws.Range("M122").Value = -2890
$ws.Range("K122").Value = 5340
$ws.Range("H122").Value = 1686.4445
$ws.Range("I122").Value = 1780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K55").Value = 6750
$ws.Range("H55").Value = 4316.25
$ws.Range("I55").Value = 2250
$ws.Range("M55").Value = -6573
$ws.Range("I132").Value = 1294
$ws.Range("J132").Value = 1075
$ws.Range("K132").Value = 11646
$ws.Range("H132").Value = 1239.25
$ws.Range("N132").Value = -14735
$ws.Range("L132").Value = 9675
$ws.Range("M132").Value = -9116
$ws.Range("M137").Value = -1341
$ws.Range("K137").Value = 6441
$ws.Range("I137").Value = 2147
$ws.Range("H137").Value = 2147

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J3").Value = 3358999.8
$ws.Range("N3").Value = -3359231.8
$ws.Range("L3").Value = 3358999.8
$ws.Range("H3").Value = 5588143
$ws.Range("H46").Value = 16000
$ws.Range("H102").Value = 3849.5
$ws.Range("I102").Value = 3819.6
$ws.Range("M102").Value = -2197.6
$ws.Range("K102").Value = 3819.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J9").Value = 3112
$ws.Range("L9").Value = 3112
$ws.Range("N9").Value = -3560
$ws.Range("H9").Value = 2305.75
$ws.Range("H30").Value = 1479.2
$ws.Range("I30").Value = 1479.2
$ws.Range("L30").Value = 0
$ws.Range("K30").Value = 1479.2
$ws.Range("M30").Value = -1371.2
$ws.Range("N30").ClearContents()
$ws.Range("J30").Value = 0
$ws.Range("H82").Value = 17368.846
$ws.Range("I82").Value = 2493.25
$ws.Range("J82").Value = 23980.223
$ws.Range("K82").Value = 2493.25
$ws.Range("M82").Value = -2132.25
$ws.Range("L82").Value = 23980.223
$ws.Range("N82").Value = -24702.223
$ws.Range("I85").Value = 2493.25
$ws.Range("M85").Value = -1245.25
$ws.Range("N85").Value = -26476.223
$ws.Range("H85").Value = 17368.846
$ws.Range("K85").Value = 2493.25
$ws.Range("J85").Value = 23980.223
$ws.Range("L85").Value = 23980.223
$ws.Range("J106").Value = 1456905.9
$ws.Range("N106").Value = -1459429.9
$ws.Range("H106").Value = 1456905.9
$ws.Range("L106").Value = 1456905.9
$ws.Range("I132").Value = 1586
$ws.Range("K132").Value = 4758
$ws.Range("H132").Value = 2196.0667
$ws.Range("M132").Value = -2228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 26010.666
$ws.Range("M20").Value = -17770
$ws.Range("K20").Value = 18010
$ws.Range("I20").Value = 18010
$ws.Range("J31").Value = 59019
$ws.Range("H31").Value = 46018.332
$ws.Range("L31").Value = 59019
$ws.Range("N31").Value = -59715
$ws.Range("J101").Value = 61357.285
$ws.Range("H101").Value = 61357.285
$ws.Range("L101").Value = 61357.285
$ws.Range("N101").Value = -67847.285
$ws.Range("L103").Value = 99763.336
$ws.Range("H103").Value = 99763.336
$ws.Range("J103").Value = 99763.336
$ws.Range("N103").Value = -102107.336
$ws.Range("H104").Value = 58689.25
$ws.Range("J104").Value = 58689.25
$ws.Range("L104").Value = 58689.25
$ws.Range("N104").Value = -65677.25
$ws.Range("I136").Value = 1399.1852
$ws.Range("M136").Value = -1647.5556
$ws.Range("H136").Value = 2108.077
$ws.Range("K136").Value = 4197.5556
